$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.212.38"
$ws.Range("E2").Value = "  -4.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.883.25"
$ws.Range("E3").Value = "  -5.47%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.29"
$ws.Range("E5").Value = "  -6.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.04"
$ws.Range("E6").Value = "  -7.59%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.415"
$ws.Range("E8").Value = "  -6.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.08"
$ws.Range("E9").Value = "  -5.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  -8.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.345"
$ws.Range("E11").Value = "  -6.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.378.61"
$ws.Range("E12").Value = "  -5.46%  "

$ws.Range("E13").Value = "  -4.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.27"
$ws.Range("E14").Value = "  -5.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000156"
$ws.Range("E15").Value = "  -8.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "55.265.37"
$ws.Range("E16").Value = "  -4.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.94"
$ws.Range("E17").Value = "  -4.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.882.28"
$ws.Range("E18").Value = "  -5.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.29"
$ws.Range("E19").Value = "  -5.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -6.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.54"
$ws.Range("E21").Value = "  -8.14%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.476"
$ws.Range("E24").Value = "  -5.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.84"
$ws.Range("E25").Value = "  -4.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.158"
$ws.Range("E27").Value = "  -6.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0832"
$ws.Range("E28").Value = "  -12.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.29"
$ws.Range("E29").Value = "  -8.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("E30").Value = "  -8.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.72"
$ws.Range("E31").Value = "  -6.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.53"
$ws.Range("E32").Value = "  -7.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.11"
$ws.Range("E33").Value = "  -10.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.23"
$ws.Range("E34").Value = "  -4.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.33"
$ws.Range("E35").Value = "  -9.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.50"
$ws.Range("E36").Value = "  -7.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.44"
$ws.Range("E37").Value = "  -3.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.17"
$ws.Range("E38").Value = "  -9.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0644"
$ws.Range("E39").Value = "  -6.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.04"
$ws.Range("E41").Value = "  -4.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.627"
$ws.Range("E42").Value = "  -5.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("E43").Value = "  -7.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.075.69"
$ws.Range("E44").Value = "  -10.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.31"
$ws.Range("E45").Value = "  -9.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.78"
$ws.Range("E46").Value = "  -4.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.896"
$ws.Range("E47").Value = "  -11.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0228"
$ws.Range("E48").Value = "  -5.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.45"
$ws.Range("E49").Value = "  -6.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0834"
$ws.Range("E50").Value = "  -7.27%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.170"
$ws.Range("E51").Value = "  -7.78%  "
